# RPA datasets push 2024-07-03
# Insert a new IPO record (이노스페이스) as the new row 2, pushing the
# existing rows (old row 2..20) down to rows 3..21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts rows 2-20 -> 3-21).
$ws.Rows.Item(2).Insert()

# The inserted row picks up the header row's bold/border formatting by
# default; clear that so the new data row matches the plain style used by
# every other data row.
$ws.Rows.Item(2).ClearFormats()

# Populate the new row with the 이노스페이스 (InnoSpace) IPO record.
$ws.Range("A2").Value = "2024-06-20"
$ws.Range("B2").Value = "이노스페이스"
$ws.Range("C2").Value = "미래"
$ws.Range("D2").Value = "2024-06-25"
$ws.Range("E2").Value = "2024-07-02"
$ws.Range("F2").Value = 57589000
$ws.Range("G2").Value = 1330000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 36400
$ws.Range("J2").Value = 43300
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 43300
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "1150.72 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"
